$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.971.32'
$ws.Range('E2').Value = '  -3.84%  '
$ws.Range('D3').Value = '3.511.88'
$ws.Range('E3').Value = '  -3.29%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '576.25'
$ws.Range('E5').Value = '  -1.95%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '170.56'
$ws.Range('E6').Value = '  -5.55%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.616'
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').Value = '3.505.27'
$ws.Range('E8').Value = '  -3.12%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.188'
$ws.Range('E10').Value = '  -6.77%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.65'
$ws.Range('E11').Value = '  +12.78%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.596'
$ws.Range('E12').Value = '  -1.79%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '47.05'
$ws.Range('E13').Value = '  -5.36%  '
$ws.Range('E14').Value = '  -3.89%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '683.49'
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('D16').Value = '4.073.09'
$ws.Range('E16').Value = '  -3.48%  '
$ws.Range('E17').Value = '  -3.23%  '
$ws.Range('D18').Value = '69.076.20'
$ws.Range('E18').Value = '  -3.87%  '
$ws.Range('D19').Value = '3.516.22'
$ws.Range('E19').Value = '  -3.83%  '
$ws.Range('E20').Value = '  -1.59%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.36'
$ws.Range('E21').Value = '  -5.09%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '11.08'
$ws.Range('E22').Value = '  -4.57%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.907'
$ws.Range('E23').Value = '  -3.36%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '16.50'
$ws.Range('E24').Value = '  -7.09%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '97.41'
$ws.Range('E25').Value = '  -5.68%  '
$ws.Range('E26').Value = '  -4.81%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.65'
$ws.Range('E28').Value = '  -6.86%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.40'
$ws.Range('E29').Value = '  -5.84%  '
$ws.Range('E30').Value = '  -5.53%  '
$ws.Range('E31').Value = '  -4.45%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.16'
$ws.Range('E32').Value = '  -7.20%  '
$ws.Range('B33').Value = 'Mantle'
$ws.Range('C33').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.35'
$ws.Range('E33').Value = '  -6.09%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '7.23'
$ws.Range('E34').Value = '  -1.48%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '578.95'
$ws.Range('E35').Value = '  -1.63%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.65'
$ws.Range('E36').Value = '  -12.65%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '10.81'
$ws.Range('E37').Value = '  -4.47%  '
$ws.Range('E38').Value = '  -3.88%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '57.36'
$ws.Range('E39').Value = '  -3.44%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.137'
$ws.Range('E41').Value = '  -3.58%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '3.452.21'
$ws.Range('E42').Value = '  -6.28%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0438'
$ws.Range('E43').Value = '  -6.58%  '
$ws.Range('E44').Value = '  -3.04%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '33.13'
$ws.Range('E45').Value = '  -7.04%  '
$ws.Range('D46').Value = '0.0₃0699'
$ws.Range('E46').Value = '  -8.36%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.88'
$ws.Range('E47').Value = '  +2.80%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.57'
$ws.Range('E48').Value = '  -7.42%  '
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '134.18'
$ws.Range('E50').Value = '  +2.06%  '
$ws.Range('E51').Value = '  -1.47%  '
